# Update "想去人数" (wanted-to-go count) figures in column F
# for both the "展览" sheet and the combined "全部类型" sheet.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - rows 2..12, column F
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 567
$wsExpo.Range("F3").Value = 186
$wsExpo.Range("F4").Value = 352
$wsExpo.Range("F5").Value = 412
$wsExpo.Range("F6").Value = 261
$wsExpo.Range("F7").Value = 2395
$wsExpo.Range("F8").Value = 411
$wsExpo.Range("F9").Value = 6196
$wsExpo.Range("F10").Value = 158
$wsExpo.Range("F11").Value = 399
$wsExpo.Range("F12").Value = 20

# Sheet "全部类型" (All types) - mirrors the exhibition rows, but rows 7,
# 8 and 14 belong to interleaved "演出" (Shows) data and stay untouched.
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 567
$wsAll.Range("F3").Value = 186
$wsAll.Range("F4").Value = 352
$wsAll.Range("F5").Value = 412
$wsAll.Range("F6").Value = 261
$wsAll.Range("F9").Value = 2395
$wsAll.Range("F10").Value = 411
$wsAll.Range("F11").Value = 6196
$wsAll.Range("F12").Value = 158
$wsAll.Range("F13").Value = 399
$wsAll.Range("F15").Value = 20
